# fix(gui) step 1 and 2
# - Bump the quote/list date in A1 by one day (45308 -> 45309).
# - Update the unit prices in column D for rows 27-34 (TEX-03..TEX-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: date bump
$ws.Range("A1").Value = 45309

# Step 2: price updates
$ws.Range("D27").Value = 13037.21
$ws.Range("D28").Value = 16298.154
$ws.Range("D29").Value = 13968.909
$ws.Range("D30").Value = 17695.7
$ws.Range("D31").Value = 14900.605
$ws.Range("D32").Value = 19089.968
$ws.Range("D33").Value = 16298.154
$ws.Range("D34").Value = 21894.902
